$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 723143
$ws.Range("I62").Value = 1145848.8
$ws.Range("J62").Value = 89084.336
$ws.Range("K62").Value = 1145848.8
$ws.Range("L62").Value = 89084.336
$ws.Range("M62").Value = -1145224.8
$ws.Range("N62").Value = -90332.336
$ws.Range("H65").Value = 723143
$ws.Range("I65").Value = 1145848.8
$ws.Range("J65").Value = 89084.336
$ws.Range("K65").Value = 5729244
$ws.Range("L65").Value = 445421.68
$ws.Range("M65").Value = -5726124
$ws.Range("N65").Value = -451661.68
$ws.Range("H92").Value = 95005.14
$ws.Range("I92").Value = 328.33334
$ws.Range("J92").Value = 221240.89
$ws.Range("K92").Value = 328.33334
$ws.Range("L92").Value = 221240.89
$ws.Range("M92").Value = 919.66666
$ws.Range("N92").Value = -223736.89
$ws.Range("H97").Value = 499.5
$ws.Range("J97").Value = 499
$ws.Range("L97").Value = 1497
$ws.Range("N97").Value = -2489
$ws.Range("H98").Value = 674.10345
$ws.Range("I98").Value = 678.3571
$ws.Range("J98").Value = 555
$ws.Range("K98").Value = 678.3571
$ws.Range("L98").Value = 555
$ws.Range("M98").Value = 819.6429
$ws.Range("N98").Value = -3551
$ws.Range("H112").Value = 2148.2964
$ws.Range("I112").Value = 7000
$ws.Range("J112").Value = 1961.6923
$ws.Range("K112").Value = 21000
$ws.Range("L112").Value = 5885.0769
$ws.Range("M112").Value = -19892
$ws.Range("N112").Value = -8101.0769
$ws.Range("H122").Value = 674.10345
$ws.Range("I122").Value = 678.3571
$ws.Range("J122").Value = 555
$ws.Range("K122").Value = 2035.0713
$ws.Range("L122").Value = 1665
$ws.Range("M122").Value = 414.9287000000002
$ws.Range("N122").Value = -6565
$ws.Range("H125").Value = 1328.8334
$ws.Range("I125").Value = 1050.25
$ws.Range("K125").Value = 9452.25
$ws.Range("M125").Value = -6992.25
$ws.Range("H131").Value = 9250.23
$ws.Range("I131").Value = 2958.8333
$ws.Range("K131").Value = 8876.499899999999
$ws.Range("M131").Value = -3836.499899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1521.3334
$ws.Range("I61").Value = 1314.2222
$ws.Range("J61").Value = 2142.6667
$ws.Range("K61").Value = 1314.2222
$ws.Range("L61").Value = 2142.6667
$ws.Range("M61").Value = -1102.2222
$ws.Range("N61").Value = -2566.6667
$ws.Range("H122").Value = 3040.862
$ws.Range("I122").Value = 1783.1428
$ws.Range("J122").Value = 4214.7334
$ws.Range("K122").Value = 5349.428400000001
$ws.Range("L122").Value = 12644.2002
$ws.Range("M122").Value = -2899.428400000001
$ws.Range("N122").Value = -17544.2002
$ws.Range("H136").Value = 1521.3334
$ws.Range("I136").Value = 1314.2222
$ws.Range("J136").Value = 2142.6667
$ws.Range("K136").Value = 3942.6666
$ws.Range("L136").Value = 6428.000100000001
$ws.Range("M136").Value = -1392.6666
$ws.Range("N136").Value = -11528.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1712.375
$ws.Range("I105").Value = 1452.9412
$ws.Range("K105").Value = 1452.9412
$ws.Range("M105").Value = 294.0588
$ws.Range("H132").Value = 98778.89
$ws.Range("J132").Value = 98778.89
$ws.Range("L132").Value = 98778.89
$ws.Range("N132").Value = -108898.89
$ws.Range("H134").Value = 1249.7
$ws.Range("I134").Value = 1110.7778
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 3332.3334
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -797.3334000000004
$ws.Range("N134").Value = -12570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2892.1538
$ws.Range("I16").Value = 1699.75
$ws.Range("J16").Value = 4800
$ws.Range("K16").Value = 1699.75
$ws.Range("L16").Value = 4800
$ws.Range("M16").Value = -1412.75
$ws.Range("N16").Value = -5374
$ws.Range("H99").Value = 2405.1875
$ws.Range("J99").Value = 2748.8333
$ws.Range("L99").Value = 2748.8333
$ws.Range("N99").Value = -5744.8333
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = ""
$ws.Range("H113").Value = 2892.1538
$ws.Range("I113").Value = 1699.75
$ws.Range("J113").Value = 4800
$ws.Range("K113").Value = 1699.75
$ws.Range("L113").Value = 4800
$ws.Range("M113").Value = 470.25
$ws.Range("N113").Value = -9140
$ws.Range("H122").Value = 2089.5264
$ws.Range("I122").Value = 2130
$ws.Range("J122").Value = 1959.1111
$ws.Range("K122").Value = 6390
$ws.Range("L122").Value = 5877.3333
$ws.Range("M122").Value = -3940
$ws.Range("N122").Value = -10777.3333
$ws.Range("H126").Value = 2405.1875
$ws.Range("J126").Value = 2748.8333
$ws.Range("L126").Value = 8246.499899999999
$ws.Range("N126").Value = -13186.4999
$ws.Range("H132").Value = 3524.1785
$ws.Range("I132").Value = 3120.95
$ws.Range("K132").Value = 9362.849999999999
$ws.Range("M132").Value = -6832.849999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 38024.93
$ws.Range("I11").Value = 54921.473
$ws.Range("J11").Value = 2354.4443
$ws.Range("K11").Value = 164764.419
$ws.Range("L11").Value = 7063.3329
$ws.Range("M11").Value = -164624.419
$ws.Range("N11").Value = -7343.3329
$ws.Range("H134").Value = 1582.25
$ws.Range("I134").Value = 1582.25
$ws.Range("K134").Value = 4746.75
$ws.Range("M134").Value = 323.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 429.25
$ws.Range("I13").Value = 502.5
$ws.Range("J13").Value = 356
$ws.Range("K13").Value = 502.5
$ws.Range("L13").Value = 356
$ws.Range("M13").Value = -363.5
$ws.Range("N13").Value = -634
$ws.Range("H70").Value = 9932.692
$ws.Range("J70").Value = 7909.3335
$ws.Range("L70").Value = 7909.3335
$ws.Range("N70").Value = -8449.3335
$ws.Range("H73").Value = 9932.692
$ws.Range("J73").Value = 7909.3335
$ws.Range("L73").Value = 7909.3335
$ws.Range("N73").Value = -9781.3335
$ws.Range("H97").Value = 1302.6129
$ws.Range("I97").Value = 1331.08
$ws.Range("J97").Value = 1184
$ws.Range("K97").Value = 1331.08
$ws.Range("L97").Value = 1184
$ws.Range("M97").Value = -835.0799999999999
$ws.Range("N97").Value = -2176
$ws.Range("H122").Value = 2076.2144
$ws.Range("I122").Value = 1761.7
$ws.Range("J122").Value = 2862.5
$ws.Range("K122").Value = 5285.1
$ws.Range("L122").Value = 8587.5
$ws.Range("M122").Value = -2835.1
$ws.Range("N122").Value = -13487.5
$ws.Range("H132").Value = 6905.5884
$ws.Range("I132").Value = 5841.6
$ws.Range("K132").Value = 17524.8
$ws.Range("M132").Value = -14994.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 22729896
$ws.Range("I7").Value = 31251732
$ws.Range("J7").Value = 5000.8335
$ws.Range("K7").Value = 31251732
$ws.Range("L7").Value = 5000.8335
$ws.Range("M7").Value = -31251620
$ws.Range("N7").Value = -5224.8335
$ws.Range("H40").Value = 3105.2222
$ws.Range("I40").Value = 2445.182
$ws.Range("K40").Value = 2445.182
$ws.Range("M40").Value = -2309.182
$ws.Range("H122").Value = 3637.0605
$ws.Range("I122").Value = 2462.423
$ws.Range("K122").Value = 7387.268999999999
$ws.Range("M122").Value = -4937.268999999999
$ws.Range("H126").Value = 22729896
$ws.Range("I126").Value = 31251732
$ws.Range("J126").Value = 5000.8335
$ws.Range("K126").Value = 93755196
$ws.Range("L126").Value = 15002.5005
$ws.Range("M126").Value = -93752726
$ws.Range("N126").Value = -19942.5005
$ws.Range("H136").Value = 1597.4
$ws.Range("I136").Value = 1436.7906
$ws.Range("K136").Value = 4310.3718
$ws.Range("M136").Value = -1760.3718

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = ""
$ws.Range("H11").Value = 15005000
$ws.Range("I11").Value = 15005000
$ws.Range("K11").Value = 15005000
$ws.Range("M11").Value = -15004858
$ws.Range("H74").Value = 28822.334
$ws.Range("J74").Value = 31266.8
$ws.Range("L74").Value = 31266.8
$ws.Range("N74").Value = -33138.8
$ws.Range("H77").Value = 28822.334
$ws.Range("J77").Value = 31266.8
$ws.Range("L77").Value = 93800.4
$ws.Range("N77").Value = -103160.4
$ws.Range("H122").Value = 1353.2632
$ws.Range("I122").Value = 1182.0333
$ws.Range("J122").Value = 1995.375
$ws.Range("K122").Value = 3546.0999
$ws.Range("L122").Value = 5986.125
$ws.Range("M122").Value = -1096.0999
$ws.Range("N122").Value = -10886.125
$ws.Range("H124").Value = 20000
$ws.Range("J124").Value = 20000
$ws.Range("L124").Value = 20000
$ws.Range("N124").Value = -29820
$ws.Range("H125").Value = 59667.43
$ws.Range("J125").Value = 62434.4
$ws.Range("L125").Value = 62434.4
$ws.Range("N125").Value = -72274.4
$ws.Range("H132").Value = 3196.9565
$ws.Range("I132").Value = 3116.0857
$ws.Range("J132").Value = 3454.2727
$ws.Range("K132").Value = 9348.2571
$ws.Range("L132").Value = 10362.8181
$ws.Range("M132").Value = -6818.257100000001
$ws.Range("N132").Value = -15422.8181
$ws.Range("H136").Value = 3204.8572
$ws.Range("I136").Value = 1225.0834
$ws.Range("K136").Value = 3675.2502
$ws.Range("M136").Value = -1125.2502
